$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on numeric-looking price cells so Excel keeps them as text (matching original inline-string semantics)
$textCells = @("D5", "D6", "D7", "D9", "D10", "D11", "D12", "D15", "D19", "D21", "D22", "D24", "D27", "D30", "D31", "D33", "D34", "D38", "D39", "D42", "D43", "D44", "D45", "D46", "D49", "D50", "D51")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply updated values
$ws.Range("D2").Value = '45.617.39'
$ws.Range("E2").Value = '  +6.86%  '
$ws.Range("D3").Value = '2.398.63'
$ws.Range("E3").Value = '  +4.13%  '
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("D5").Value = '115.29'
$ws.Range("E5").Value = '  +10.39%  '
$ws.Range("D6").Value = '319.77'
$ws.Range("E6").Value = '  +3.10%  '
$ws.Range("D7").Value = '0.637'
$ws.Range("E7").Value = '  +2.99%  '
$ws.Range("E8").Value = '  -0.20%  '
$ws.Range("D9").Value = '0.629'
$ws.Range("E9").Value = '  +3.89%  '
$ws.Range("D10").Value = '42.74'
$ws.Range("E10").Value = '  +8.04%  '
$ws.Range("D11").Value = '0.0932'
$ws.Range("E11").Value = '  +3.23%  '
$ws.Range("D12").Value = '8.72'
$ws.Range("E12").Value = '  +5.41%  '
$ws.Range("E13").Value = '  +3.33%  '
$ws.Range("E14").Value = '  +2.86%  '
$ws.Range("D15").Value = '16.02'
$ws.Range("E15").Value = '  +4.44%  '
$ws.Range("D16").Value = '2.764.38'
$ws.Range("E16").Value = '  -0.70%  '
$ws.Range("D17").Value = '2.398.26'
$ws.Range("E17").Value = '  +4.35%  '
$ws.Range("D18").Value = '45.624.08'
$ws.Range("E18").Value = '  +6.54%  '
$ws.Range("D19").Value = '7.52'
$ws.Range("E19").Value = '  +2.75%  '
$ws.Range("E20").Value = '  +3.76%  '
$ws.Range("D21").Value = '13.68'
$ws.Range("E21").Value = '  +1.90%  '
$ws.Range("D22").Value = '75.05'
$ws.Range("E22").Value = '  +2.30%  '
$ws.Range("E23").Value = '  +4.49%  '
$ws.Range("D24").Value = '265.16'
$ws.Range("E24").Value = '  -1.14%  '
$ws.Range("E25").Value = '  +7.01%  '
$ws.Range("E26").Value = '  -0.66%  '
$ws.Range("D27").Value = '7.76'
$ws.Range("E27").Value = '  +6.08%  '
$ws.Range("E29").Value = '  +2.58%  '
$ws.Range("D30").Value = '40.08'
$ws.Range("E30").Value = '  +10.54%  '
$ws.Range("D31").Value = '0.0993'
$ws.Range("E31").Value = '  +15.69%  '
$ws.Range("E32").Value = '  +2.60%  '
$ws.Range("D33").Value = '173.22'
$ws.Range("E33").Value = '  +5.04%  '
$ws.Range("D34").Value = '2.94'
$ws.Range("E34").Value = '  +11.86%  '
$ws.Range("E35").Value = '  +2.16%  '
$ws.Range("E36").Value = '  +11.09%  '
$ws.Range("E37").Value = '  +6.94%  '
$ws.Range("D38").Value = '4.20'
$ws.Range("E38").Value = '  +15.80%  '
$ws.Range("D39").Value = '3.10'
$ws.Range("E39").Value = '  +10.13%  '
$ws.Range("E40").Value = '  +5.19%  '
$ws.Range("E41").Value = '  +12.76%  '
$ws.Range("B42").Value = 'Algorand'
$ws.Range("C42").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D42").Value = '0.243'
$ws.Range("E42").Value = '  +7.44%  '
$ws.Range("B43").Value = 'Celestia'
$ws.Range("C43").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D43").Value = '13.78'
$ws.Range("E43").Value = '  +11.90%  '
$ws.Range("D44").Value = '100.49'
$ws.Range("E44").Value = '  -8.20%  '
$ws.Range("D45").Value = '72.26'
$ws.Range("E45").Value = '  +2.03%  '
$ws.Range("D46").Value = '90.01'
$ws.Range("E46").Value = '  +15.61%  '
$ws.Range("E47").Value = '  -0.44%  '
$ws.Range("E48").Value = '  +14.29%  '
$ws.Range("D49").Value = '116.46'
$ws.Range("E49").Value = '  +5.09%  '
$ws.Range("D50").Value = '9.50'
$ws.Range("E50").Value = '  +10.05%  '
$ws.Range("D51").Value = '1.59'
$ws.Range("E51").Value = '  +10.62%  '
